$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the redundant second header row (units row: '‰air', '‰ V-PDB', 'ratio').
# This shifts all data rows up by one, and Excel will drop the now-unused
# shared strings / styles for that row automatically on save.
$ws.Rows(2).Delete()

# Update the visible selection/cursor to match the post-edit state.
$ws.Range("B16").Select()
